$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("T" tonnage) values were changed from plain numbers to text
# values with a trailing "T" suffix (e.g. 45 -> "45T"). G2 originally used
# a quote-prefixed text style, so we re-apply it with a leading apostrophe
# so the cell keeps its original (quotePrefix) style instead of drifting to
# a plain-number style.
$ws.Range("G2").Value = "'45T"
$ws.Range("G3").Value = "300T"
$ws.Range("G4").Value = "400T"
$ws.Range("G5").Value = "200T"
$ws.Range("G6").Value = "100T"
$ws.Range("G7").Value = "45T"
$ws.Range("G8").Value = "90T"
$ws.Range("G9").Value = "130T"
$ws.Range("G10").Value = "90T"
$ws.Range("G11").Value = "500T"
$ws.Range("G12").Value = "90T"
$ws.Range("G13").Value = "160T"
$ws.Range("G14").Value = "130T"
$ws.Range("G15").Value = "300T"
$ws.Range("G16").Value = "160T"
$ws.Range("G17").Value = "500T"
$ws.Range("G18").Value = "250T"
$ws.Range("G19").Value = "180T"
$ws.Range("G20").Value = "180T"
$ws.Range("G21").Value = "250T"
$ws.Range("G22").Value = "100T"
$ws.Range("G23").Value = "160T"
$ws.Range("G24").Value = "160T"
$ws.Range("G25").Value = "160T"
$ws.Range("G26").Value = "160T"
$ws.Range("G27").Value = "160T"
$ws.Range("G28").Value = "200T"
$ws.Range("G29").Value = "750T"
$ws.Range("G30").Value = "100T"
$ws.Range("G31").Value = "80T"
$ws.Range("G32").Value = "80T"
$ws.Range("G33").Value = "80T"
$ws.Range("G34").Value = "80T"
$ws.Range("G35").Value = "600T"
$ws.Range("G36").Value = "80T"
$ws.Range("G37").Value = "80T"
$ws.Range("G38").Value = "80T"
$ws.Range("G39").Value = "80T"
$ws.Range("G40").Value = "100T"
$ws.Range("G41").Value = "90T"
$ws.Range("G42").Value = "90T"
$ws.Range("G43").Value = "130T"
$ws.Range("G44").Value = "160T"
$ws.Range("G45").Value = "160T"
$ws.Range("G46").Value = "160T"
$ws.Range("G47").Value = "80T"
$ws.Range("G48").Value = "80T"
$ws.Range("G49").Value = "80T"
$ws.Range("G50").Value = "80T"
$ws.Range("G51").Value = "80T"
$ws.Range("G52").Value = "1250T"
$ws.Range("G53").Value = "110T"
$ws.Range("G54").Value = "110T"
$ws.Range("G55").Value = "110T"
$ws.Range("G56").Value = "110T"
$ws.Range("G57").Value = "110T"
$ws.Range("G58").Value = "110T"
$ws.Range("G59").Value = "110T"
$ws.Range("G60").Value = "160T"
$ws.Range("G61").Value = "160T"
$ws.Range("G62").Value = "160T"
$ws.Range("G63").Value = "160T"
$ws.Range("G64").Value = "160T"
$ws.Range("G65").Value = "800T"
$ws.Range("G66").Value = "650T"
$ws.Range("G67").Value = "650T"
$ws.Range("G69").Value = "336T"
$ws.Range("G75").Value = "336T"
$ws.Range("G147").Value = "800T"
$ws.Range("G152").Value = "350T"
$ws.Range("G153").Value = "60T"
$ws.Range("G154").Value = "60T"
$ws.Range("G155").Value = "60T"
$ws.Range("G164").Value = "160T"
$ws.Range("G165").Value = "160T"
$ws.Range("G169").Value = "3T"
$ws.Range("G170").Value = "3T"
$ws.Range("G171").Value = "10T"

# Rename the internal Data Model linked-table defined name (cosmetic
# artifact of Excel's Data Model sync after the range's data changed).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlcn.WorksheetConnection_VehiclesothersA1K1501") {
        $n.Name = "_xlcn.WorksheetConnection_VehiclesothersA1K150"
    }
}

# Move the active selection to G7, matching the edited workbook's view state.
$ws.Range("G7").Select()
